$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 130, pushing existing rows 130-242 down to 131-243.
$ws.Rows(130).Insert()

# Populate the newly inserted row 130 with the new record's data.
$ws.Cells.Item(130, 1).Value = 10
$ws.Cells.Item(130, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(130, 3).Value = "La Araucanía"
$ws.Cells.Item(130, 4).Value = 44827
$ws.Cells.Item(130, 5).Value = 9
$ws.Cells.Item(130, 6).Value = 100112005
$ws.Cells.Item(130, 7).Value = "Puerro"
$ws.Cells.Item(130, 8).Value = "Azul de Maquehue"
$ws.Cells.Item(130, 9).Value = "Primera"
$ws.Cells.Item(130, 10).Value = 20
$ws.Cells.Item(130, 11).Value = 15000
$ws.Cells.Item(130, 12).Value = 15000
$ws.Cells.Item(130, 13).Value = 15000
$ws.Cells.Item(130, 14).Value = "$/docena de paquetes"
$ws.Cells.Item(130, 15).Value = "Provincia de Cautín"
$ws.Cells.Item(130, 16).Value = 1250
$ws.Cells.Item(130, 17).Value = 12
$ws.Cells.Item(130, 18).Value = "Hortaliza"
